$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2021-10-11"

# Update the row label for October to reflect the new "through" date
$ws.Range("A11").Value = "October (through 10-11)"

# Update the October row (row 11) values
$ws.Range("B11").Value = 10
$ws.Range("C11").Value = 18
$ws.Range("D11").Value = 18
$ws.Range("E11").Value = 29
$ws.Range("G11").Value = 50
$ws.Range("H11").Value = 72

# Update the Total row (row 12) values
$ws.Range("B12").Value = 236
$ws.Range("C12").Value = 447
$ws.Range("D12").Value = 645
$ws.Range("E12").Value = 577
$ws.Range("G12").Value = 951
$ws.Range("H12").Value = 1322
